# Code restructure for SearchAuthorClusterPage and minor change in WAT03 testscript
#
# - WAT03 (row 4) JIRA ID cell (B4) changes from "OPQA-TBD" to the list of
#   JIRA tickets now associated with the restructured Author search page.
# - The active selection on the sheet moves from E4 to B4 (the cell that was
#   just edited).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the JIRA ID for WAT03 (row 4, column B).
$ws.Range("B4").Value = "WAT-162||WAT-507||WAT-215||WAT-220"

# Move / leave the selection on the cell that was edited.
$ws.Range("B4").Select()
